$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MLM_Summary")

# Row 2 - Logistic_Regression
$ws.Range("B2").Value = 0.87
$ws.Range("C2").Value = 0.84
$ws.Range("E2").Value = 39
$ws.Range("F2").Value = 6
$ws.Range("H2").Value = 0.68
$ws.Range("J2").Value = 0.72
$ws.Range("K2").Value = 0.87
$ws.Range("L2").Value = 0.82

# Row 3 - Random_Forest
$ws.Range("B3").Value = 0.95
$ws.Range("C3").Value = 0.82
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 0.6899999999999999
$ws.Range("I3").Value = 0.65
$ws.Range("K3").Value = 0.89
$ws.Range("L3").Value = 0.77

# Row 4 - Kernel_SVM
$ws.Range("B4").Value = 0.9
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 38
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 0.68
$ws.Range("I4").Value = 0.88
$ws.Range("J4").Value = 0.77
$ws.Range("K4").Value = 0.84
$ws.Range("L4").Value = 0.86

# Row 5 - CatBoost
$ws.Range("B5").Value = 0.96
$ws.Range("C5").Value = 0.84
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 41
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 0.73
$ws.Range("I5").Value = 0.65
$ws.Range("J5").Value = 0.6899999999999999
$ws.Range("K5").Value = 0.91
$ws.Range("L5").Value = 0.78

# Row 6 - DNN
$ws.Range("B6").Value = 0.98
$ws.Range("C6").Value = 0.8100000000000001
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 12
$ws.Range("H6").Value = 0.5600000000000001
$ws.Range("J6").Value = 0.68
$ws.Range("K6").Value = 0.73
$ws.Range("L6").Value = 0.8100000000000001
